# "44344 consultation proposal line space"
#
# The paragraph holding the "<Proposal Description>" placeholder (the
# one wrapped by the _GoBack bookmark) is split in two:
#   1. a new, empty paragraph using the same paragraph formatting
#      (Keybody style, italic paragraph mark) is inserted immediately
#      before it -- this is what gives the placeholder its own blank
#      line of space above it;
#   2. the placeholder keeps its Keybody/italic formatting, and the
#      bookmark (now collapsed -- start immediately followed by end)
#      plus its run stay together on the second paragraph.

$d = $word.ActiveDocument

$bm = $d.Bookmarks.Item("_GoBack")

# Capture the existing placeholder text before we rewrite the range.
$placeholderText = $bm.Range.Text

# Grab the whole paragraph (including its end-of-paragraph mark) so the
# replacement XML can cleanly introduce a second paragraph break.
$pRange = $bm.Range.Duplicate
$pRange.Expand(4)  # wdParagraph

$escapedText = $placeholderText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

$xml = '<?xml version="1.0" standalone="yes"?>' + `
'<?mso-application progid="Word.Document"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
          '<w:p>' + `
            '<w:pPr>' + `
              '<w:pStyle w:val="Keybody"/>' + `
              '<w:rPr><w:i/><w:iCs/></w:rPr>' + `
            '</w:pPr>' + `
          '</w:p>' + `
          '<w:p>' + `
            '<w:pPr>' + `
              '<w:pStyle w:val="Keybody"/>' + `
              '<w:rPr><w:i/><w:iCs/></w:rPr>' + `
            '</w:pPr>' + `
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
            '<w:bookmarkEnd w:id="0"/>' + `
            '<w:r>' + `
              '<w:rPr><w:i/><w:iCs/></w:rPr>' + `
              "<w:t>$escapedText</w:t>" + `
            '</w:r>' + `
          '</w:p>' + `
        '</w:body>' + `
      '</w:document>' + `
    '</pkg:xmlData>' + `
  '</pkg:part>' + `
'</pkg:package>'

$pRange.InsertXML($xml)
